# Fix logging system configuration
# Appends a new data row (row 75) to each of the four log sheets,
# mirroring the existing row-74 formatting (date/time style for column A,
# plain text for columns B-E, plain numbers for F-I).

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"
$newTime = 45861.46377314815

$rowsData = @{
    "MID_LFT_#1" = @{
        B = "0x01,0x90"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
        D = "0x01,0x40"
        E = "0x07"
        F = 400
        G = "5.68631262647113e+23"
        H = 320
        I = 7
    }
    "MID_LFT_#2" = @{
        B = "0x01,0x7c"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
        D = "0x01,0x3C"
        E = "0x19"
        F = 380
        G = "5.68432987514711e+23"
        H = 316
        I = 25
    }
    "MID_PLT_#1" = @{
        B = "0x00,0x6e"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
        D = "0x00,0x62"
        E = "0x15"
        F = 110
        G = "5.68631262647113e+23"
        H = 98
        I = 15
    }
    "MID_PLT_#2" = @{
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
        D = "0x00,0x77"
        E = "0x9"
        F = 130
        G = "5.68631262647113e+23"
        H = 119
        I = 9
    }
}

$sheetNames = @("MID_LFT_#1", "MID_LFT_#2", "MID_PLT_#1", "MID_PLT_#2")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $data = $rowsData[$name]

    $ws.Range("A75").Value() = $newTime
    $ws.Range("A75").NumberFormat = $dateFormat

    $ws.Range("B75").Value() = $data.B
    $ws.Range("C75").Value() = $data.C
    $ws.Range("D75").Value() = $data.D
    $ws.Range("E75").Value() = $data.E

    $ws.Range("F75").Value() = $data.F
    $g = [double]$data.G
    $ws.Range("G75").Value() = $g
    $ws.Range("H75").Value() = $data.H
    $ws.Range("I75").Value() = $data.I
}
